$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H126").Value2 = 0
$ws.Range("J126").Value2 = 0
$ws.Range("L126").Value2 = 0
$ws.Range("N126").ClearContents()

$ws.Range("H127").Value2 = 2712.5557
$ws.Range("J127").Value2 = 4216
$ws.Range("L127").Value2 = 12648
$ws.Range("N127").Value2 = -22568

$ws.Range("H129").Value2 = 7989.375
$ws.Range("I129").Value2 = 995.6667
$ws.Range("J129").Value2 = 12185.6
$ws.Range("K129").Value2 = 2987.0001
$ws.Range("L129").Value2 = 36556.8
$ws.Range("M129").Value2 = 2012.9999
$ws.Range("N129").Value2 = -46556.8

$ws.Range("H130").Value2 = 0
$ws.Range("J130").Value2 = 0
$ws.Range("L130").Value2 = 0
$ws.Range("N130").ClearContents()

$ws.Range("H132").Value2 = 2043.75
$ws.Range("I132").Value2 = 2091.6667
$ws.Range("K132").Value2 = 6275.000100000001
$ws.Range("M132").Value2 = -3745.000100000001

$ws.Range("H137").Value2 = 421524.7
$ws.Range("I137").Value2 = 3090.3635
$ws.Range("K137").Value2 = 9271.0905
$ws.Range("M137").Value2 = -6721.0905

$ws.Range("H138").Value2 = 2069.0862
$ws.Range("I138").Value2 = 712.43243
$ws.Range("J138").Value2 = 4459.381
$ws.Range("K138").Value2 = 2137.29729
$ws.Range("L138").Value2 = 13378.143
$ws.Range("M138").Value2 = 3002.70271
$ws.Range("N138").Value2 = -23658.143

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value2 = 8480.574000000001
$ws.Range("I32").Value2 = 5752.5386
$ws.Range("J32").Value2 = 24242.555
$ws.Range("K32").Value2 = 5752.5386
$ws.Range("L32").Value2 = 24242.555
$ws.Range("M32").Value2 = -5465.5386
$ws.Range("N32").Value2 = -24816.555

$ws.Range("H74").Value2 = 40465.414
$ws.Range("I74").Value2 = 22202.209
$ws.Range("J74").Value2 = 302238
$ws.Range("K74").Value2 = 22202.209
$ws.Range("L74").Value2 = 302238
$ws.Range("M74").Value2 = -21328.209
$ws.Range("N74").Value2 = -303986

$ws.Range("H77").Value2 = 40465.414
$ws.Range("I77").Value2 = 22202.209
$ws.Range("J77").Value2 = 302238
$ws.Range("K77").Value2 = 111011.045
$ws.Range("L77").Value2 = 1511190
$ws.Range("M77").Value2 = -106643.045
$ws.Range("N77").Value2 = -1519926

$ws.Range("H86").Value2 = 115156.5
$ws.Range("J86").Value2 = 115156.5
$ws.Range("L86").Value2 = 115156.5
$ws.Range("N86").Value2 = -117528.5

$ws.Range("H89").Value2 = 115156.5
$ws.Range("J89").Value2 = 115156.5
$ws.Range("L89").Value2 = 345469.5
$ws.Range("N89").Value2 = -357325.5

$ws.Range("H122").Value2 = 27233.057
$ws.Range("I122").Value2 = 1230.9032
$ws.Range("K122").Value2 = 3692.7096
$ws.Range("M122").Value2 = -1242.7096

$ws.Range("H132").Value2 = 2212.372
$ws.Range("I132").Value2 = 1937.3206
$ws.Range("K132").Value2 = 5811.9618
$ws.Range("M132").Value2 = -3281.9618

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H45").Value2 = 19000
$ws.Range("J45").Value2 = 19000
$ws.Range("L45").Value2 = 19000
$ws.Range("N45").Value2 = -20616

$ws.Range("H94").Value2 = 150798.17
$ws.Range("I94").Value2 = 0
$ws.Range("K94").Value2 = 0
$ws.Range("M94").ClearContents()

$ws.Range("H97").Value2 = 5667
$ws.Range("I97").Value2 = 1000.5
$ws.Range("J97").Value2 = 15000
$ws.Range("K97").Value2 = 1000.5
$ws.Range("L97").Value2 = 15000
$ws.Range("M97").Value2 = -9.5
$ws.Range("N97").Value2 = -16982

$ws.Range("H105").Value2 = 3147.1304
$ws.Range("I105").Value2 = 2995.111
$ws.Range("K105").Value2 = 2995.111
$ws.Range("M105").Value2 = -1248.111

$ws.Range("H107").Value2 = 9355.790999999999
$ws.Range("I107").Value2 = 9050
$ws.Range("K107").Value2 = 9050
$ws.Range("M107").Value2 = -7130

$ws.Range("H134").Value2 = 2924.2295
$ws.Range("I134").Value2 = 2429.6738
$ws.Range("J134").Value2 = 4440.8667
$ws.Range("K134").Value2 = 7289.0214
$ws.Range("L134").Value2 = 13322.6001
$ws.Range("M134").Value2 = -4754.0214
$ws.Range("N134").Value2 = -18392.6001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value2 = 2752.147
$ws.Range("I31").Value2 = 2529.037
$ws.Range("J31").Value2 = 2899.0732
$ws.Range("K31").Value2 = 2529.037
$ws.Range("L31").Value2 = 2899.0732
$ws.Range("M31").Value2 = -2234.037
$ws.Range("N31").Value2 = -3489.0732

$ws.Range("H34").Value2 = 2752.147
$ws.Range("I34").Value2 = 2529.037
$ws.Range("J34").Value2 = 2899.0732
$ws.Range("K34").Value2 = 2529.037
$ws.Range("L34").Value2 = 2899.0732
$ws.Range("M34").Value2 = -2327.037
$ws.Range("N34").Value2 = -3303.0732

$ws.Range("H58").Value2 = 8050.12
$ws.Range("I58").Value2 = 3824.611
$ws.Range("J58").Value2 = 18915.715
$ws.Range("K58").Value2 = 3824.611
$ws.Range("L58").Value2 = 18915.715
$ws.Range("M58").Value2 = -3621.611
$ws.Range("N58").Value2 = -19321.715

$ws.Range("H88").Value2 = 28493.428
$ws.Range("J88").Value2 = 27857.166
$ws.Range("L88").Value2 = 27857.166
$ws.Range("N88").Value2 = -28669.166

$ws.Range("H91").Value2 = 28493.428
$ws.Range("J91").Value2 = 27857.166
$ws.Range("L91").Value2 = 27857.166
$ws.Range("N91").Value2 = -30665.166

$ws.Range("H99").Value2 = 3378.625
$ws.Range("I99").Value2 = 3602.8572
$ws.Range("J99").Value2 = 2950.5454
$ws.Range("K99").Value2 = 3602.8572
$ws.Range("L99").Value2 = 2950.5454
$ws.Range("M99").Value2 = -2104.8572
$ws.Range("N99").Value2 = -5946.5454

$ws.Range("H126").Value2 = 3378.625
$ws.Range("I126").Value2 = 3602.8572
$ws.Range("J126").Value2 = 2950.5454
$ws.Range("K126").Value2 = 10808.5716
$ws.Range("L126").Value2 = 8851.636200000001
$ws.Range("M126").Value2 = -8338.571599999999
$ws.Range("N126").Value2 = -13791.6362

$ws.Range("H132").Value2 = 3731.1462
$ws.Range("I132").Value2 = 1910.5217
$ws.Range("K132").Value2 = 5731.5651
$ws.Range("M132").Value2 = -3201.5651

$ws.Range("H134").Value2 = 2510.6667
$ws.Range("I134").Value2 = 2249.0393
$ws.Range("J134").Value2 = 6958.3335
$ws.Range("K134").Value2 = 6747.117899999999
$ws.Range("L134").Value2 = 20875.0005
$ws.Range("M134").Value2 = -4212.117899999999
$ws.Range("N134").Value2 = -25945.0005

$ws.Range("H136").Value2 = 8050.12
$ws.Range("I136").Value2 = 3824.611
$ws.Range("J136").Value2 = 18915.715
$ws.Range("K136").Value2 = 11473.833
$ws.Range("L136").Value2 = 56747.145
$ws.Range("M136").Value2 = -8923.832999999999
$ws.Range("N136").Value2 = -61847.145

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value2 = 150220.83
$ws.Range("J122").Value2 = 450115
$ws.Range("L122").Value2 = 4051035
$ws.Range("N122").Value2 = -4055935

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value2 = 4185.5527
$ws.Range("I132").Value2 = 4174.5864
$ws.Range("J132").Value2 = 4220.8887
$ws.Range("K132").Value2 = 12523.7592
$ws.Range("L132").Value2 = 12662.6661
$ws.Range("M132").Value2 = -9993.7592
$ws.Range("N132").Value2 = -17722.6661

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value2 = 7382.927
$ws.Range("I22").Value2 = 8215.333000000001
$ws.Range("K22").Value2 = 8215.333000000001
$ws.Range("M22").Value2 = -7920.333000000001

$ws.Range("H27").Value2 = 7382.927
$ws.Range("I27").Value2 = 8215.333000000001
$ws.Range("K27").Value2 = 8215.333000000001
$ws.Range("M27").Value2 = -8108.333000000001

$ws.Range("H58").Value2 = 9845.5
$ws.Range("I58").Value2 = 6243.25
$ws.Range("K58").Value2 = 6243.25
$ws.Range("M58").Value2 = -5983.25

$ws.Range("H69").Value2 = 23081.5
$ws.Range("J69").Value2 = 23081.5
$ws.Range("L69").Value2 = 23081.5
$ws.Range("N69").Value2 = -24703.5

$ws.Range("H72").Value2 = 23081.5
$ws.Range("J72").Value2 = 23081.5
$ws.Range("L72").Value2 = 69244.5
$ws.Range("N72").Value2 = -77356.5

$ws.Range("H93").Value2 = 1293.4546
$ws.Range("I93").Value2 = 1567.3334
$ws.Range("J93").Value2 = 706.5714
$ws.Range("K93").Value2 = 1567.3334
$ws.Range("L93").Value2 = 706.5714
$ws.Range("M93").Value2 = -319.3334
$ws.Range("N93").Value2 = -3202.5714

$ws.Range("H135").Value2 = 0
$ws.Range("J135").Value2 = 0
$ws.Range("L135").Value2 = 0
$ws.Range("N135").ClearContents()

$ws.Range("H136").Value2 = 27418
$ws.Range("I136").Value2 = 2162.9355
$ws.Range("K136").Value2 = 6488.806500000001
$ws.Range("M136").Value2 = -3938.806500000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H92").Value2 = 116662.336
$ws.Range("J92").Value2 = 116662.336
$ws.Range("L92").Value2 = 116662.336
$ws.Range("N92").Value2 = -121654.336

$ws.Range("H132").Value2 = 2698.0527
$ws.Range("I132").Value2 = 2561.2727
$ws.Range("J132").Value2 = 3600.8
$ws.Range("K132").Value2 = 7683.8181
$ws.Range("L132").Value2 = 10802.4
$ws.Range("M132").Value2 = -5153.8181
$ws.Range("N132").Value2 = -15862.4

$ws.Range("H136").Value2 = 3305.2615
$ws.Range("I136").Value2 = 3081.07
$ws.Range("K136").Value2 = 9243.210000000001
$ws.Range("M136").Value2 = -6693.210000000001
